$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's data started on row 2 (row 1 was blank). Deleting row 1
# shifts the whole data block up by one row, matching the diff.
$ws.Rows("1:1").Delete()
